$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set "andet" in column E for rows 24-31 (Wednesday column, same value as C/D columns)
foreach ($r in 24..31) {
    $ws.Range("E$r").Value = "andet"
}

# Update selection to E32 (matches diff's sheetView selection change)
$ws.Range("E32").Select()
